# Rename the three header/footer picture placeholders:
#   footer (first-page)  : image2.png -> image1.png   (id=3)
#   footer (primary)     : image2.png -> image1.png   (id=2)
#   header (first-page)  : image1.jpg -> image2.jpg   (id=1)
#
# InlineShape.Name is not writable directly in this object model, so we
# round-trip each inline picture through ConvertToShape()/Name=.../
# ConvertToInlineShape() to rename it while keeping it inline.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-InlinePicture($range, $newName) {
    $inline = $range.InlineShapes.Item(1)
    $shape = $inline.ConvertToShape()
    $shape.Name = $newName
    $shape.ConvertToInlineShape() | Out-Null
}

# Footer (first page) -- Pearson logo "image2.png" -> "image1.png"
Rename-InlinePicture $sec.Footers.Item(2).Range "image1.png"

# Footer (primary/default) -- Pearson logo "image2.png" -> "image1.png"
Rename-InlinePicture $sec.Footers.Item(1).Range "image1.png"

# Header (first page) -- BTEC logo "image1.jpg" -> "image2.jpg"
Rename-InlinePicture $sec.Headers.Item(2).Range "image2.jpg"
